# Adds a new "Tema 10 - Adjectives 2 - Adjectives to describe" vocabulary
# block to the "Lista vocabulario" workbook:
#   - 40 new English/Spanish adjective rows appended to Hoja1 (rows 267-306)
#   - a new "Hoja2" sheet that mirrors the English/Spanish pairs in a
#     two-column "English" / "Español" layout (rows 7-46)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The 40 new adjective word-pairs (English, Spanish), in entry order.
$pairs = @(
    @("Spacious", "Espacioso"),
    @("Comfortable", "Cómodo"),
    @("Cozy", "Acogedor"),
    @("Bright", "Brillante / Luminoso"),
    @("Modern", "Moderno"),
    @("Traditional", "Tradicional"),
    @("Elegant", "Elegante"),
    @("Simple", "Sencillo"),
    @("Quiet", "Tranquilo"),
    @("Noisy", "Ruidoso"),
    @("Tidy", "Ordenado"),
    @("Messy", "Desordenado"),
    @("Big", "Grande"),
    @("Small", "Pequeño"),
    @("Clean", "Limpio"),
    @("Dark", "Oscuro"),
    @("Warm", "Cálido"),
    @("Cold", "Frío"),
    @("New", "Nuevo"),
    @("Old", "Viejo"),
    @("Soft", "Suave"),
    @("Hard", "Duro"),
    @("Neat", "Ordenado"),
    @("Wide", "Ancho"),
    @("Narrow", "Estrecho"),
    @("Beautiful", "Hermoso"),
    @("Organised", "Organizado"),
    @("Colourful", "Colorido"),
    @("Attractive", "Atractivo"),
    @("Peaceful", "Pacífico"),
    @("Decorative", "Decorativo"),
    @("Welcoming", "Acogedor"),
    @("Convenient", "Conveniente"),
    @("Fantastic", "Fantástico"),
    @("Impressive", "Impresionante"),
    @("Brightened", "Iluminado"),
    @("Appealing", "Atractivo"),
    @("Refreshing", "Refrescante"),
    @("Wonderful", "Maravilloso"),
    @("Artistic", "Artístico")
)

$topicShort = "Adjectives 2 - Adjectives to describe"
$topicFull = "Tema 10 - Adjectives 2 - Adjectives to describe"

# --- Hoja1: append the new rows right after the existing data (row 266) ---
# Column order mirrors the original authoring session: first the English /
# Spanish word pairs for every new row, then the topic name (F267), then the
# "Tema" column (E267:E306), and finally the concatenation formula (H267).
$firstRow = 267
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $r = $firstRow + $i
    $en = $pairs[$i][0]
    $es = $pairs[$i][1]

    $ws1.Cells.Item($r, 1).Value = $en   # A: Inglés
    $ws1.Cells.Item($r, 3).Value = $es   # C: Español
}

# Only the first new row also carries the "Nombre tema" (F) value.
$ws1.Cells.Item($firstRow, 6).Value = $topicShort

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $r = $firstRow + $i
    $ws1.Cells.Item($r, 5).Value = $topicFull  # E: Tema
}

$ws1.Range("H267").Formula = '=CONCATENATE(E267," ","-"," ",F267)'

# --- Hoja2: new sheet placed right after Hoja1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

$ws2.Cells.Item(6, 5).Value = "English"
$ws2.Cells.Item(6, 6).Value = "Español"

$headerRow = 7
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $r = $headerRow + $i
    $ws2.Cells.Item($r, 5).Value = $pairs[$i][0]  # E: English
    $ws2.Cells.Item($r, 6).Value = $pairs[$i][1]  # F: Español
}

$ws2.Range("E7:F46").Select()

# Leave Hoja1 as the active sheet/selection when done.
$ws1.Activate()
$ws1.Range("E267:E306").Select()
